$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all affected cells to Text format first so values are stored as literal strings
# (matching the original inlineStr/text representation), then assign new values,
# then reset the style back to Normal so no extra style index is introduced.

$changedRefs = @("D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5", "E5", "G5", "D6", "E6", "G6", "D7", "E7", "G7", "D8", "E8", "G8", "D9", "E9", "G9", "D10", "E10", "G10", "D11", "E11", "G11", "D12", "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "D15", "E15", "G15", "D16", "E16", "G16", "D17", "E17", "G17", "D18", "E18", "G18", "E19", "G19", "E20", "G20", "D21", "E21", "G21", "D22", "E22", "G22", "D23", "G23", "E24", "G24", "D25", "E25", "G25", "E26", "G26", "E27", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35", "G36", "G37", "G38", "G39", "D40", "E40", "G40", "D41", "E41", "G41", "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "E45", "G45", "D46", "E46", "G46", "D47", "E47", "G47", "G48", "G49", "G50", "G51")
foreach ($ref in $changedRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "275.63"
$ws.Range("E2").Value = "0.95%"
$ws.Range("G2").Value = "3"
$ws.Range("D3").Value = "27.18"
$ws.Range("E3").Value = "1.85%"
$ws.Range("G3").Value = "3"
$ws.Range("D4").Value = "4.859"
$ws.Range("E4").Value = "-0.81%"
$ws.Range("G4").Value = "3"
$ws.Range("D5").Value = "0.06403"
$ws.Range("E5").Value = "1.11%"
$ws.Range("G5").Value = "3"
$ws.Range("D6").Value = "6.926"
$ws.Range("E6").Value = "0.25%"
$ws.Range("G6").Value = "3"
$ws.Range("D7").Value = "1.223"
$ws.Range("E7").Value = "-3.44%"
$ws.Range("G7").Value = "3"
$ws.Range("D8").Value = "0.8806"
$ws.Range("E8").Value = "0.28%"
$ws.Range("G8").Value = "3"
$ws.Range("D9").Value = "0.1513"
$ws.Range("E9").Value = "3.84%"
$ws.Range("G9").Value = "3"
$ws.Range("D10").Value = "0.05089"
$ws.Range("E10").Value = "2.69%"
$ws.Range("G10").Value = "3"
$ws.Range("D11").Value = "0.07617"
$ws.Range("E11").Value = "3.73%"
$ws.Range("G11").Value = "3"
$ws.Range("D12").Value = "0.02988"
$ws.Range("E12").Value = "-4.40%"
$ws.Range("G12").Value = "3"
$ws.Range("D13").Value = "0.08997"
$ws.Range("E13").Value = "-0.55%"
$ws.Range("G13").Value = "3"
$ws.Range("D14").Value = "0.001567"
$ws.Range("E14").Value = "-0.47%"
$ws.Range("G14").Value = "3"
$ws.Range("D15").Value = "0.0006414"
$ws.Range("E15").Value = "1.17%"
$ws.Range("G15").Value = "3"
$ws.Range("D16").Value = "0.006201"
$ws.Range("E16").Value = "2.72%"
$ws.Range("G16").Value = "3"
$ws.Range("D17").Value = "3.464"
$ws.Range("E17").Value = "0.26%"
$ws.Range("G17").Value = "3"
$ws.Range("D18").Value = "3.309"
$ws.Range("E18").Value = "-1.51%"
$ws.Range("G18").Value = "3"
$ws.Range("E19").Value = "0.54%"
$ws.Range("G19").Value = "3"
$ws.Range("E20").Value = "-0.95%"
$ws.Range("G20").Value = "3"
$ws.Range("D21").Value = "0.1355"
$ws.Range("E21").Value = "2.07%"
$ws.Range("G21").Value = "3"
$ws.Range("D22").Value = "3.920"
$ws.Range("E22").Value = "0.14%"
$ws.Range("G22").Value = "3"
$ws.Range("D23").Value = "0.04427"
$ws.Range("G23").Value = "3"
$ws.Range("E24").Value = "-0.29%"
$ws.Range("G24").Value = "3"
$ws.Range("D25").Value = "0.004269"
$ws.Range("E25").Value = "15.55%"
$ws.Range("G25").Value = "3"
$ws.Range("E26").Value = "-0.13%"
$ws.Range("G26").Value = "3"
$ws.Range("E27").Value = "13.75%"
$ws.Range("G27").Value = "3"
$ws.Range("G28").Value = "3"
$ws.Range("G29").Value = "3"
$ws.Range("G30").Value = "3"
$ws.Range("G31").Value = "3"
$ws.Range("G32").Value = "3"
$ws.Range("G33").Value = "3"
$ws.Range("G34").Value = "3"
$ws.Range("G35").Value = "3"
$ws.Range("G36").Value = "3"
$ws.Range("G37").Value = "3"
$ws.Range("G38").Value = "3"
$ws.Range("G39").Value = "3"
$ws.Range("D40").Value = "0.04138"
$ws.Range("E40").Value = "2.46%"
$ws.Range("G40").Value = "3"
$ws.Range("D41").Value = "0.006831"
$ws.Range("E41").Value = "2.45%"
$ws.Range("G41").Value = "3"
$ws.Range("E42").Value = "0.67%"
$ws.Range("G42").Value = "3"
$ws.Range("D43").Value = "0.002154"
$ws.Range("E43").Value = "2.25%"
$ws.Range("G43").Value = "3"
$ws.Range("D44").Value = "0.01186"
$ws.Range("E44").Value = "0.54%"
$ws.Range("G44").Value = "3"
$ws.Range("E45").Value = "-2.85%"
$ws.Range("G45").Value = "3"
$ws.Range("D46").Value = "1.653"
$ws.Range("E46").Value = "-47.10%"
$ws.Range("G46").Value = "3"
$ws.Range("D47").Value = "0.02003"
$ws.Range("E47").Value = "-0.08%"
$ws.Range("G47").Value = "3"
$ws.Range("G48").Value = "3"
$ws.Range("G49").Value = "3"
$ws.Range("G50").Value = "3"
$ws.Range("G51").Value = "3"

foreach ($ref in $changedRefs) {
    $ws.Range($ref).Style = "Normal"
}

